$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row below mirrors one <row> block changed in the source diff.
# D-column (Price) values are forced to Text via a leading apostrophe,
# then the style is reset to "Normal" so no stray number-format style
# sticks to the cell (the workbook's Price column is text, e.g. '222.63').

# Row 2
$ws.Range('D2').Value = "'34.530.20"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.53%  '

# Row 3
$ws.Range('D3').Value = "'1.788.06"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.03%  '

# Row 4
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').Value = "'222.64"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.57%  '

# Row 6
$ws.Range('E6').Value = '  -1.25%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').Value = "'32.13"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.94%  '

# Row 9
$ws.Range('E9').Value = '  +0.17%  '

# Row 10
$ws.Range('E10').Value = '  +2.81%  '

# Row 11
$ws.Range('E11').Value = '  +1.29%  '

# Row 12
$ws.Range('D12').Value = "'2.044.50"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.17%  '

# Row 13
$ws.Range('D13').Value = "'1.790.63"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.00%  '

# Row 14
$ws.Range('D14').Value = "'10.93"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.68%  '

# Row 15
$ws.Range('D15').Value = "'34.503.26"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.19%  '

# Row 16
$ws.Range('E16').Value = '  +0.71%  '

# Row 17
$ws.Range('E17').Value = '  +1.90%  '

# Row 18
$ws.Range('D18').Value = "'68.69"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.48%  '

# Row 19
$ws.Range('D19').Value = "'253.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.82%  '

# Row 20
$ws.Range('D20').Value = "'0.0₃0783"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.59%  '

# Row 21
$ws.Range('E21').Value = '  +0.11%  '

# Row 22
$ws.Range('D22').Value = "'10.46"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.08%  '

# Row 23
$ws.Range('D23').Value = "'4.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.53%  '

# Row 24
$ws.Range('E24').Value = '  +0.33%  '

# Row 25
$ws.Range('D25').Value = "'160.55"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.31%  '

# Row 26
$ws.Range('D26').Value = "'16.38"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.79%  '

# Row 27
$ws.Range('D27').Value = "'7.08"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.16%  '

# Row 28
$ws.Range('E28').Value = '  -0.87%  '

# Row 29
$ws.Range('E29').Value = '  +0.02%  '

# Row 30
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'3.75"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.96%  '

# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = "'0.0516"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.19%  '

# Row 32
$ws.Range('E32').Value = '  -0.55%  '

# Row 33
$ws.Range('E33').Value = '  -0.79%  '

# Row 34
$ws.Range('D34').Value = "'1.88"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.87%  '

# Row 35
$ws.Range('D35').Value = "'1.436.98"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.45%  '

# Row 36
$ws.Range('D36').Value = "'0.638"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.50%  '

# Row 37
$ws.Range('E37').Value = '  -1.43%  '

# Row 38
$ws.Range('E38').Value = '  +2.27%  '

# Row 39
$ws.Range('D39').Value = "'84.68"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.41%  '

# Row 40
$ws.Range('D40').Value = "'2.80"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.74%  '

# Row 41
$ws.Range('D41').Value = "'2.34"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.42%  '

# Row 42
$ws.Range('D42').Value = "'0.916"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.51%  '

# Row 43
$ws.Range('D43').Value = "'2.08"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.46%  '

# Row 44
$ws.Range('D44').Value = "'6.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.51%  '

# Row 45
$ws.Range('E45').Value = '  -1.09%  '

# Row 46
$ws.Range('E46').Value = '  -4.69%  '

# Row 47
$ws.Range('D47').Value = "'1.944.70"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.07%  '

# Row 48
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'12.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.84%  '

# Row 49
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'1.00"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.07%  '

# Row 50
$ws.Range('D50').Value = "'103.76"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.80%  '

# Row 51
$ws.Range('D51').Value = "'49.80"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.18%  '
